# Case_4_34 parallel.xlsx update: extend table with two more columns (P, Q)
# and fix up the contingency values in columns I, K, M, O for rows 2-25.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (header row): add P1=14 and Q1=15, matching O1's header style ---
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# --- Data rows 2-25: flip I/K and M/O values, and append new P/Q columns ---
for ($r = 2; $r -le 25; $r++) {
    $ws.Range("I$r").Value = 2
    $ws.Range("K$r").Value = 1
    $ws.Range("M$r").Value = 2
    $ws.Range("O$r").Value = 1
    $ws.Range("P$r").Value = 2
    $ws.Range("Q$r").Value = 2
}
